$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells keep their text formatting (values like "536.29" would
# otherwise be auto-converted to numbers by Excel's type inference).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.415.01"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "3.153.61"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "536.29"
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("D6").Value = "139.62"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +8.91%  "
$ws.Range("D9").Value = "7.35"
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("E11").Value = "  +5.27%  "
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").Value = "3.696.72"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").Value = "25.71"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("E15").Value = "  +6.69%  "
$ws.Range("D16").Value = "58.474.88"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "3.150.25"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "6.20"
$ws.Range("E18").Value = "  +6.20%  "
$ws.Range("D19").Value = "13.02"
$ws.Range("E19").Value = "  +5.26%  "
$ws.Range("E20").Value = "  +6.04%  "
$ws.Range("D21").Value = "372.93"
$ws.Range("E21").Value = "  +7.40%  "
$ws.Range("D22").Value = "5.79"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "70.05"
$ws.Range("D25").Value = "0.515"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "8.14"
$ws.Range("E28").Value = "  +14.66%  "
$ws.Range("D29").Value = "0.0₃0868"
$ws.Range("E29").Value = "  +3.42%  "
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("E31").Value = "  +5.71%  "
$ws.Range("D32").Value = "21.91"
$ws.Range("E32").Value = "  +4.73%  "
$ws.Range("D33").Value = "5.17"
$ws.Range("E33").Value = "  +8.14%  "
$ws.Range("E34").Value = "  +5.28%  "
$ws.Range("D35").Value = "160.58"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").Value = "6.23"
$ws.Range("E36").Value = "  +4.62%  "
$ws.Range("D37").Value = "1.37"
$ws.Range("E37").Value = "  +13.53%  "
$ws.Range("D38").Value = "25.32"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "2.641.46"
$ws.Range("E39").Value = "  +9.75%  "
$ws.Range("E40").Value = "  +6.69%  "
$ws.Range("E41").Value = "  +4.34%  "
$ws.Range("E42").Value = "  +4.48%  "
$ws.Range("D43").Value = "38.70"
$ws.Range("E43").Value = "  +5.99%  "
$ws.Range("D44").Value = "0.707"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("E45").Value = "  +8.98%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "3.197.45"
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("E48").Value = "  +10.97%  "
$ws.Range("D49").Value = "6.21"
$ws.Range("E49").Value = "  +4.37%  "
$ws.Range("D50").Value = "0.978"
$ws.Range("E50").Value = "  +5.18%  "
$ws.Range("E51").Value = "  +4.92%  "
